$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 738 (shifts existing rows 738+ down to 740+)
$ws.Range("A738:A739").EntireRow.Insert()

# Row 738: 2026/02/01, Sun, 22, 21
# Force column A to be read as literal text so the date-like string "2026/02/01"
# isn't auto-converted into a date serial number, then restore the Normal
# style so no stray number-format style lingers on the cell.
$ws.Range("A738").NumberFormat = "@"
$ws.Range("A738").Value = "2026/02/01"
$ws.Range("A738").Style = "Normal"
$ws.Range("B738").Value = "日"
$ws.Range("C738").Value = 22
$ws.Range("D738").Value = 21

# Row 739: 2026/02/02, Mon, 1, 22
$ws.Range("A739").NumberFormat = "@"
$ws.Range("A739").Value = "2026/02/02"
$ws.Range("A739").Style = "Normal"
$ws.Range("B739").Value = "月"
$ws.Range("C739").Value = 1
$ws.Range("D739").Value = 22
